# Update cryptos list values (D: Price, E: Volume(1h)) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to text format so values like "237.40" or
# "0.06195" are not auto-converted to numbers by Excel on assignment.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '25.798.54'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '1.756.95'
$ws.Range('E3').Value = '  +1.03%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '237.40'
$ws.Range('E5').Value = '  +0.62%  '
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').Value = '0.5064'
$ws.Range('E7').Value = '  +3.47%  '
$ws.Range('D8').Value = '41.07'
$ws.Range('E8').Value = '  -2.17%  '
$ws.Range('E9').Value = '  +8.68%  '
$ws.Range('D10').Value = '0.06195'
$ws.Range('E10').Value = '  +2.97%  '
$ws.Range('D11').Value = '1.754.48'
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('D12').Value = '0.06954'
$ws.Range('E12').Value = '  +4.46%  '
$ws.Range('D13').Value = '15.45'
$ws.Range('E13').Value = '  +6.47%  '
$ws.Range('D14').Value = '0.6013'
$ws.Range('E14').Value = '  +2.30%  '
$ws.Range('D15').Value = '77.61'
$ws.Range('E15').Value = '  +0.38%  '
$ws.Range('D16').Value = '4.448'
$ws.Range('E16').Value = '  +2.56%  '
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('E18').Value = '  -0.06%  '
$ws.Range('D19').Value = '25.831.04'
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').Value = '11.64'
$ws.Range('E20').Value = '  +3.35%  '
$ws.Range('D21').Value = '0.000006793'
$ws.Range('E21').Value = '  +8.18%  '
$ws.Range('D22').Value = '1.978.25'
$ws.Range('E22').Value = '  +0.65%  '
$ws.Range('D23').Value = '4.058'
$ws.Range('E23').Value = '  +4.63%  '
$ws.Range('D24').Value = '8.141'
$ws.Range('E24').Value = '  +2.54%  '
$ws.Range('D25').Value = '5.169'
$ws.Range('E25').Value = '  +1.44%  '
$ws.Range('D26').Value = '137.88'
$ws.Range('E26').Value = '  +1.75%  '
$ws.Range('D27').Value = '1.465'
$ws.Range('E27').Value = '  -0.95%  '
$ws.Range('D28').Value = '15.01'
$ws.Range('E28').Value = '  +4.97%  '
$ws.Range('D29').Value = '1.811'
$ws.Range('E29').Value = '  -2.40%  '
$ws.Range('D30').Value = '102.60'
$ws.Range('E30').Value = '  +3.50%  '
$ws.Range('D31').Value = '0.08264'
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('D32').Value = '3.681'
$ws.Range('E32').Value = '  +1.42%  '
$ws.Range('D33').Value = '3.389'
$ws.Range('E33').Value = '  +6.42%  '
$ws.Range('D34').Value = '0.04374'
$ws.Range('E34').Value = '  +1.85%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Value = '2.655'
$ws.Range('E36').Value = '  +1.55%  '
$ws.Range('D37').Value = '1.003'
$ws.Range('E37').Value = '  -2.68%  '
$ws.Range('D38').Value = '0.6009'
$ws.Range('E38').Value = '  -0.54%  '
$ws.Range('D39').Value = '2.726'
$ws.Range('E39').Value = '  -3.48%  '
$ws.Range('D40').Value = '1.949'
$ws.Range('E40').Value = '  -5.57%  '
$ws.Range('D41').Value = '0.01546'
$ws.Range('E41').Value = '  +3.73%  '
$ws.Range('D42').Value = '1.002'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').Value = '103.40'
$ws.Range('E43').Value = '  +2.35%  '
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('D45').Value = '0.7440'
$ws.Range('E45').Value = '  -5.59%  '
$ws.Range('D46').Value = '4.857'
$ws.Range('E46').Value = '  -5.30%  '
$ws.Range('D47').Value = '0.05478'
$ws.Range('E47').Value = '  +8.03%  '
$ws.Range('D48').Value = '0.1076'
$ws.Range('E48').Value = '  +4.14%  '
$ws.Range('D49').Value = '5.941'
$ws.Range('E49').Value = '  -1.94%  '
$ws.Range('D50').Value = '30.19'
$ws.Range('E50').Value = '  +2.33%  '
$ws.Range('D51').Value = '1.001'
$ws.Range('E51').Value = '  +0.19%  '

# Restore the default (Normal) style on column D so no stray number-format
# style index is left applied to the cells.
$priceRange.Style = "Normal"

